$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Classes")
$v = $ws.Range("E2").Value
Write-Host "E2 type:" $v.GetType().FullName
Write-Host "E2 value:" $v
